# Weekly update of Fruta/Hortaliza market data rows.
# Applies new Fecha / Volumen / Precio (min, max, promedio) / Origen / Precio $/Kg
# values to rows 3-21 (row 2 and row 17 are unchanged), reflecting the latest
# weekly data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3,4).Value = 44873
$ws.Cells.Item(3,11).Value = 8000
$ws.Cells.Item(3,12).Value = 8000
$ws.Cells.Item(3,13).Value = 8000
$ws.Cells.Item(3,16).Value = 500
$ws.Cells.Item(4,4).Value = 44846
$ws.Cells.Item(4,10).Value = 250
$ws.Cells.Item(4,12).Value = 5000
$ws.Cells.Item(4,13).Value = 5000
$ws.Cells.Item(4,16).Value = 312
$ws.Cells.Item(5,4).Value = 44231
$ws.Cells.Item(5,10).Value = 250
$ws.Cells.Item(5,11).Value = 5000
$ws.Cells.Item(5,13).Value = 5500
$ws.Cells.Item(5,15).Value = "Provincia de Quillota"
$ws.Cells.Item(5,16).Value = 344
$ws.Cells.Item(6,4).Value = 44251
$ws.Cells.Item(6,10).Value = 120
$ws.Cells.Item(6,12).Value = 5000
$ws.Cells.Item(6,13).Value = 5000
$ws.Cells.Item(6,15).Value = "Región Metropolitana"
$ws.Cells.Item(6,16).Value = 312
$ws.Cells.Item(7,4).Value = 44188
$ws.Cells.Item(7,10).Value = 210
$ws.Cells.Item(8,4).Value = 44230
$ws.Cells.Item(8,12).Value = 6000
$ws.Cells.Item(8,13).Value = 5500
$ws.Cells.Item(8,16).Value = 344
$ws.Cells.Item(9,4).Value = 44204
$ws.Cells.Item(9,10).Value = 430
$ws.Cells.Item(10,4).Value = 44189
$ws.Cells.Item(10,10).Value = 250
$ws.Cells.Item(10,11).Value = 5000
$ws.Cells.Item(10,12).Value = 6000
$ws.Cells.Item(10,13).Value = 5500
$ws.Cells.Item(10,16).Value = 344
$ws.Cells.Item(11,4).Value = 44208
$ws.Cells.Item(12,4).Value = 44236
$ws.Cells.Item(12,10).Value = 180
$ws.Cells.Item(12,11).Value = 4000
$ws.Cells.Item(12,12).Value = 4500
$ws.Cells.Item(12,13).Value = 4167
$ws.Cells.Item(12,15).Value = "Región Metropolitana"
$ws.Cells.Item(12,16).Value = 260
$ws.Cells.Item(13,4).Value = 44855
$ws.Cells.Item(13,10).Value = 70
$ws.Cells.Item(13,11).Value = 6000
$ws.Cells.Item(13,12).Value = 7000
$ws.Cells.Item(13,13).Value = 6500
$ws.Cells.Item(13,16).Value = 406
$ws.Cells.Item(14,4).Value = 44210
$ws.Cells.Item(14,10).Value = 340
$ws.Cells.Item(14,11).Value = 5000
$ws.Cells.Item(14,12).Value = 6000
$ws.Cells.Item(14,13).Value = 5500
$ws.Cells.Item(14,16).Value = 344
$ws.Cells.Item(15,4).Value = 44186
$ws.Cells.Item(15,10).Value = 160
$ws.Cells.Item(16,4).Value = 44187
$ws.Cells.Item(16,10).Value = 160
$ws.Cells.Item(16,11).Value = 5000
$ws.Cells.Item(16,12).Value = 6000
$ws.Cells.Item(16,13).Value = 5500
$ws.Cells.Item(16,16).Value = 344
$ws.Cells.Item(18,4).Value = 44232
$ws.Cells.Item(18,10).Value = 250
$ws.Cells.Item(18,12).Value = 6000
$ws.Cells.Item(18,13).Value = 5500
$ws.Cells.Item(18,15).Value = "Provincia de Quillota"
$ws.Cells.Item(18,16).Value = 344
$ws.Cells.Item(19,4).Value = 44882
$ws.Cells.Item(19,10).Value = 70
$ws.Cells.Item(19,11).Value = 7000
$ws.Cells.Item(19,12).Value = 7000
$ws.Cells.Item(19,13).Value = 7000
$ws.Cells.Item(19,16).Value = 438
$ws.Cells.Item(20,4).Value = 44292
$ws.Cells.Item(20,10).Value = 90
$ws.Cells.Item(20,11).Value = 6000
$ws.Cells.Item(20,12).Value = 6000
$ws.Cells.Item(20,13).Value = 6000
$ws.Cells.Item(20,16).Value = 375
$ws.Cells.Item(21,4).Value = 44883
$ws.Cells.Item(21,10).Value = 180
$ws.Cells.Item(21,11).Value = 7000
$ws.Cells.Item(21,12).Value = 8000
$ws.Cells.Item(21,13).Value = 7500
$ws.Cells.Item(21,16).Value = 469
